$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.807.65'
$ws.Range('E2').Value = '  +4.29%  '
$ws.Range('D3').Value = '2.421.56'
$ws.Range('E3').Value = '  +2.61%  '
$c = $ws.Range('D4')
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range('E4').Value = '  +0.04%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '316.33'
$c.Style = "Normal"
$ws.Range('E5').Value = '  +4.42%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '101.81'
$c.Style = "Normal"
$ws.Range('E6').Value = '  +6.81%  '
$c = $ws.Range('D7')
$c.NumberFormat = "@"
$c.Value = '0.514'
$c.Style = "Normal"
$ws.Range('E7').Value = '  +2.49%  '
$ws.Range('E8').Value = '  -0.02%  '
$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.533'
$c.Style = "Normal"
$ws.Range('E9').Value = '  +11.96%  '
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '35.41'
$c.Style = "Normal"
$ws.Range('E10').Value = '  +3.02%  '
$ws.Range('E11').Value = '  +1.78%  '
$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '18.73'
$c.Style = "Normal"
$ws.Range('E12').Value = '  +1.66%  '
$ws.Range('E13').Value = '  -1.99%  '
$ws.Range('D15').Value = '2.800.80'
$ws.Range('E15').Value = '  +2.75%  '
$ws.Range('D16').Value = '2.408.53'
$ws.Range('E16').Value = '  +2.33%  '
$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '0.834'
$c.Style = "Normal"
$ws.Range('E17').Value = '  +4.77%  '
$ws.Range('D18').Value = '44.612.67'
$ws.Range('E18').Value = '  +3.84%  '
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '12.39'
$c.Style = "Normal"
$ws.Range('E19').Value = '  +4.01%  '
$ws.Range('E20').Value = '  +1.91%  '
$ws.Range('E21').Value = '  +3.73%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '68.76'
$c.Style = "Normal"
$ws.Range('E22').Value = '  +1.29%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '242.40'
$c.Style = "Normal"
$ws.Range('E23').Value = '  +3.10%  '
$ws.Range('E24').Value = '  +4.37%  '
$ws.Range('E25').Value = '  +2.35%  '
$ws.Range('E26').Value = '  -0.12%  '
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '25.17'
$c.Style = "Normal"
$ws.Range('E28').Value = '  -3.70%  '
$ws.Range('E29').Value = '  +1.65%  '
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '33.65'
$c.Style = "Normal"
$ws.Range('E30').Value = '  +4.23%  '
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '48.68'
$c.Style = "Normal"
$ws.Range('E31').Value = '  +1.74%  '
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '0.126'
$c.Style = "Normal"
$ws.Range('E32').Value = '  +18.32%  '
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '19.53'
$c.Style = "Normal"
$ws.Range('E33').Value = '  +11.61%  '
$ws.Range('E34').Value = '  +3.28%  '
$ws.Range('E35').Value = '  +0.26%  '
$ws.Range('E36').Value = '  +6.27%  '
$ws.Range('E37').Value = '  +3.71%  '
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '4.49'
$c.Style = "Normal"
$ws.Range('E38').Value = '  +3.74%  '
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '2.85'
$c.Style = "Normal"
$ws.Range('E39').Value = '  +0.64%  '
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '123.08'
$c.Style = "Normal"
$ws.Range('E40').Value = '  -3.71%  '
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '0.109'
$c.Style = "Normal"
$ws.Range('E41').Value = '  +1.85%  '
$ws.Range('E42').Value = '  -2.80%  '
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '20.89'
$c.Style = "Normal"
$ws.Range('E43').Value = '  +0.92%  '
$ws.Range('E44').Value = '  +4.32%  '
$ws.Range('D45').Value = '1.943.57'
$ws.Range('B46').Value = 'ApeXProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '2.11'
$c.Style = "Normal"
$ws.Range('E46').Value = '  -1.68%  '
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '2.94'
$c.Style = "Normal"
$ws.Range('E47').Value = '  +8.28%  '
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '9.29'
$c.Style = "Normal"
$ws.Range('E48').Value = '  +0.45%  '
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '1.75'
$c.Style = "Normal"
$ws.Range('E49').Value = '  +16.37%  '
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '75.71'
$c.Style = "Normal"
$ws.Range('E50').Value = '  +6.07%  '
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '53.96'
$c.Style = "Normal"
$ws.Range('E51').Value = '  +5.47%  '
